$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Match formatting of the row above (row 27) for the new row, then fill values.
$ws.Range("A27:C27").Copy()
$ws.Range("A28:C28").PasteSpecial(-4122)
$ws.Rows.Item(28).RowHeight = $ws.Rows.Item(27).RowHeight

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Windows"
$ws.Range("C28").Value = "Download Directory"

[void]$ws.Range("C28").Select()
